$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:33:55.427084",
    "2021-10-05 14:33:55.427092",
    "2021-10-05 14:33:55.427096",
    "2021-10-05 14:33:55.427099",
    "2021-10-05 14:33:55.427102",
    "2021-10-05 14:33:55.427105",
    "2021-10-05 14:33:55.427107",
    "2021-10-05 14:33:55.427110",
    "2021-10-05 14:33:55.427113",
    "2021-10-05 14:33:55.427116",
    "2021-10-05 14:33:55.427119",
    "2021-10-05 14:33:55.427121",
    "2021-10-05 14:33:55.427124",
    "2021-10-05 14:33:55.427126",
    "2021-10-05 14:33:55.427129"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $newTimes[$i]
}

# Insert the new "metadata" sheet right after the "data" sheet.
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$header = $meta.Range("B1:G1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

$meta.Range("A2").Value = 0
$meta.Range("A2").Font.Bold = $true
$meta.Range("A2").Borders.LineStyle = 1
$meta.Range("A2").HorizontalAlignment = -4108
$meta.Range("A2").VerticalAlignment = -4160
$meta.Range("B2").Value = "Haematuria_Alport"
$meta.Range("C2").Value = 39
$meta.Range("D2").Value = "'1.0"
$meta.Range("E2").Value = "2021-01-20T10:11:52.933176Z"
$meta.Range("F2").Value = "2021-10-05 14:33:55.423318"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/39/?format=json"
